$d = $word.ActiveDocument

# The document's first paragraph contains a hyperlink field whose visible display
# text "link" is stored as two runs - "li" followed by "nk" - with a _GoBack
# bookmark sitting between them (both runs share the rStyle "3" / FollowedHyperlink
# character style). The edit re-splits that same "link" text into "l" + "ink",
# i.e. the bookmark's structural position between the two runs does not move, but
# one character ("i") is reassigned from the end of the first run to the start of
# the second run.

# Scope the search to the first ~1000 characters so we only match the hyperlink's
# "link" text near the top of the document, and not the unrelated "Material link"
# text that appears much later in the file.
$scope = $d.Range(0, 1000)
$found = $scope.Find.Execute("link", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)

if ($found) {
    $linkStart = $scope.Start
    $linkEnd = $scope.End

    # First run currently holds "li" (2 chars) -> becomes "l".
    $run1 = $d.Range($linkStart, $linkStart + 2)
    $run1.Text = "l"

    # Second run currently holds "nk"; after run1 shrank by one character its
    # start position shifts left by one -> becomes "ink".
    $run2 = $d.Range($linkStart + 1, $linkEnd - 1)
    $run2.Text = "ink"
}
